# Add a mandatory PAN column for each investor / entity.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the existing "Tags" column (B) so the
# layout becomes: Name* | PAN* | Tags | Category* | City | Fund*
$null = $ws.Columns.Item(2).Insert()

# Header for the new column
$ws.Range("B1").Value = "PAN *"

# PAN values for each of the 6 investors
$ws.Range("B2").Value = "BUHNXDFEA6"
$ws.Range("B3").Value = "JN2GOV5FYI"
$ws.Range("B4").Value = "CGKT9ROWB1"
$ws.Range("B5").Value = "4I3FNDATK0"
$ws.Range("B6").Value = "5AM81UTOQB"
$ws.Range("B7").Value = "QNEL3S7Z2J"

# Keep a couple of formatted-but-empty cells below the table, mirroring the
# blank placeholder cells that already existed under the other columns.
$ws.Range("B8").Style = "Normal"
$ws.Range("B9").Style = "Normal"

# The new PAN column should be as wide as the Name column.
$ws.Columns.Item(1).ColumnWidth = 34.17
$ws.Columns.Item(2).ColumnWidth = 34.17

# Match the updated selection left behind in the saved workbook.
$null = $ws.Range("B11").Select()

Write-Host "PAN column added"
